$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.867.36'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '2.486.27'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'587.88"
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = "'0.515"
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = "'0.143"
$ws.Range("E9").Value = '  +3.13%  '
$ws.Range("E10").Value = '  -1.65%  '
$ws.Range("D11").Value = "'4.97"
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = "'0.332"
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").Value = '2.937.14'
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("D14").Value = "'25.32"
$ws.Range("E14").Value = '  -0.65%  '
$ws.Range("D15").Value = '67.777.58'
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("D17").Value = '2.464.90'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").Value = "'10.82"
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("D19").Value = "'7.35"
$ws.Range("E19").Value = '  -3.30%  '
$ws.Range("D20").Value = "'346.13"
$ws.Range("E20").Value = '  -1.03%  '
$ws.Range("D21").Value = "'4.11"
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = "'70.80"
$ws.Range("E23").Value = '  +2.71%  '
$ws.Range("D24").Value = "'4.15"
$ws.Range("E24").Value = '  -2.05%  '
$ws.Range("D25").Value = "'1.70"
$ws.Range("E25").Value = '  -6.08%  '
$ws.Range("D26").Value = "'8.86"
$ws.Range("E26").Value = '  -2.82%  '
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").Value = '0.0₃0885'
$ws.Range("E29").Value = '  -2.62%  '
$ws.Range("D30").Value = "'497.01"
$ws.Range("E30").Value = '  -1.50%  '
$ws.Range("D31").Value = "'7.74"
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = "'164.29"
$ws.Range("E35").Value = '  +1.27%  '
$ws.Range("E36").Value = '  +1.66%  '
$ws.Range("D37").Value = "'18.63"
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("D38").Value = "'18.17"
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = "'1.31"
$ws.Range("E40").Value = '  -1.93%  '
$ws.Range("D41").Value = "'1.72"
$ws.Range("E41").Value = '  +2.36%  '
$ws.Range("E42").Value = '  -1.18%  '
$ws.Range("D43").Value = "'4.75"
$ws.Range("E43").Value = '  -1.25%  '
$ws.Range("D44").Value = "'2.36"
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").Value = "'148.09"
$ws.Range("E45").Value = '  +3.55%  '
$ws.Range("D46").Value = "'3.52"
$ws.Range("E46").Value = '  +1.40%  '
$ws.Range("D47").Value = "'0.511"
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("D48").Value = '0.0₆0252'
$ws.Range("E48").Value = '  -4.95%  '
$ws.Range("D49").Value = "'0.0733"
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("E50").Value = '  -1.78%  '
$ws.Range("E51").Value = '  -1.54%  '
